# Update openpyxl command to current version
# Rename data files from the 2019-12-19 batch to the 2020-02-12 batch,
# and refresh the recomputed concentration values (CellConcentration /
# PipetteConcentration) that shifted by a floating point ULP when the
# experiment sheet was regenerated with the current openpyxl version.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("plate")

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $cell.Value2 -replace "20191219", "20200212"
}

$ws.Range("G5").Value = 4.633136414336805
$ws.Range("F6").Value = 0.02500000000000001
$ws.Range("G6").Value = 4.633136414336805
$ws.Range("G7").Value = 4.633136414336805
$ws.Range("F8").Value = 0.02500000000000001
$ws.Range("G8").Value = 4.633136414336805
$ws.Range("G9").Value = 4.633136414336805
$ws.Range("F10").Value = 0.02500000000000001
$ws.Range("G10").Value = 4.633136414336805
